$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("Förändrad") holds a date serial number (45181 = 2023-09-12) for
# every data row (2..176). Bump it by one day (45182 = 2023-09-13) to reflect
# the automatic update.
for ($row = 2; $row -le 176; $row++) {
    $ws.Cells.Item($row, 3).Value2 = 45182
}
